$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44271
$ws.Range("J2").Value = 200
$ws.Range("M2").Value = 1920
$ws.Range("P2").Value = 1920

# Row 3
$ws.Range("D3").Value = 44260
$ws.Range("J3").Value = 220
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1909
$ws.Range("N3").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("P3").Value = 1909
$ws.Range("Q3").Value = 1

# Row 4
$ws.Range("D4").Value = 44264
$ws.Range("J4").Value = 130
$ws.Range("M4").Value = 1908
$ws.Range("P4").Value = 1908

# Row 5
$ws.Range("D5").Value = 44266
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 1800
$ws.Range("M5").Value = 1913
$ws.Range("P5").Value = 1913

# Row 6
$ws.Range("D6").Value = 44524
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 2000

# Row 7
$ws.Range("D7").Value = 44265
$ws.Range("J7").Value = 220
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 1909
$ws.Range("N7").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O7").Value = 'Provincia de Diguillín'
$ws.Range("P7").Value = 1909
$ws.Range("Q7").Value = 1

# Row 8
$ws.Range("D8").Value = 44272
$ws.Range("J8").Value = 150
$ws.Range("M8").Value = 1893
$ws.Range("P8").Value = 1893

# Row 9
$ws.Range("D9").Value = 44263
$ws.Range("J9").Value = 140
$ws.Range("M9").Value = 1914
$ws.Range("P9").Value = 1914

# Row 10
$ws.Range("D10").Value = 44159
$ws.Range("J10").Value = 55
$ws.Range("K10").Value = 7000
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 7455
$ws.Range("N10").Value = '$/caja 36 atados'
$ws.Range("O10").Value = 'Región Metropolitana'
$ws.Range("P10").Value = 207
$ws.Range("Q10").Value = 36

# Row 11
$ws.Range("D11").Value = 44208
$ws.Range("J11").Value = 130
$ws.Range("M11").Value = 1908
$ws.Range("O11").Value = 'Provincia de Cautín'
$ws.Range("P11").Value = 1908

# Row 12
$ws.Range("D12").Value = 44166
$ws.Range("J12").Value = 240
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 700
$ws.Range("M12").Value = 633
$ws.Range("P12").Value = 633

# Row 13
$ws.Range("D13").Value = 44267
$ws.Range("J13").Value = 150
$ws.Range("M13").Value = 1913
$ws.Range("P13").Value = 1913

# Row 15
$ws.Range("D15").Value = 44525
$ws.Range("J15").Value = 60

# Row 16
$ws.Range("D16").Value = 44160
$ws.Range("J16").Value = 190
$ws.Range("K16").Value = 1300
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = 1395
$ws.Range("N16").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("P16").Value = 930
$ws.Range("Q16").Value = 1.5

# Row 17
$ws.Range("D17").Value = 44211
$ws.Range("J17").Value = 120
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 1883
$ws.Range("P17").Value = 1883

# Row 18
$ws.Range("D18").Value = 44273
$ws.Range("J18").Value = 140
$ws.Range("M18").Value = 1914
$ws.Range("P18").Value = 1914
